$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 5820.6665
$ws.Range("J40").Value = 2981
$ws.Range("L40").Value = 2981
$ws.Range("N40").Value = -3331
$ws.Range("H76").Value = 9025.526
$ws.Range("I76").Value = 8632.532999999999
$ws.Range("K76").Value = 8632.532999999999
$ws.Range("M76").Value = -8317.532999999999
$ws.Range("H79").Value = 9025.526
$ws.Range("I79").Value = 8632.532999999999
$ws.Range("K79").Value = 8632.532999999999
$ws.Range("M79").Value = -7540.532999999999
$ws.Range("H112").Value = 1988.6207
$ws.Range("I112").Value = 1298.25
$ws.Range("J112").Value = 2099.08
$ws.Range("K112").Value = 3894.75
$ws.Range("L112").Value = 6297.24
$ws.Range("M112").Value = -2786.75
$ws.Range("N112").Value = -8513.24
$ws.Range("H130").Value = 66250
$ws.Range("J130").Value = 66250
$ws.Range("L130").Value = 66250
$ws.Range("N130").Value = -76290
$ws.Range("H135").Value = 1698.75
$ws.Range("I135").Value = 1698.75
$ws.Range("K135").Value = 15288.75
$ws.Range("M135").Value = -12753.75
$ws.Range("H137").Value = 4522.095
$ws.Range("I137").Value = 2145.4546
$ws.Range("J137").Value = 7136.4
$ws.Range("K137").Value = 6436.3638
$ws.Range("L137").Value = 21409.2
$ws.Range("M137").Value = -3886.3638
$ws.Range("N137").Value = -26509.2
$ws.Range("H138").Value = 4672.2563
$ws.Range("J138").Value = 5028.6206
$ws.Range("L138").Value = 15085.8618
$ws.Range("N138").Value = -25365.8618
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1064.5714
$ws.Range("I45").Value = 1064.5714
$ws.Range("K45").Value = 1064.5714
$ws.Range("M45").Value = -687.5714
$ws.Range("H61").Value = 2849.0688
$ws.Range("I61").Value = 2213.3
$ws.Range("K61").Value = 2213.3
$ws.Range("M61").Value = -2001.3
$ws.Range("H63").Value = 2528.375
$ws.Range("I63").Value = 2445.1765
$ws.Range("K63").Value = 2445.1765
$ws.Range("M63").Value = -1759.1765
$ws.Range("H66").Value = 2528.375
$ws.Range("I66").Value = 2445.1765
$ws.Range("K66").Value = 12225.8825
$ws.Range("M66").Value = -8793.8825
$ws.Range("H132").Value = 4252.385
$ws.Range("I132").Value = 1658.5
$ws.Range("J132").Value = 8402.6
$ws.Range("K132").Value = 4975.5
$ws.Range("L132").Value = 25207.8
$ws.Range("M132").Value = -2445.5
$ws.Range("N132").Value = -30267.8
$ws.Range("H136").Value = 2849.0688
$ws.Range("I136").Value = 2213.3
$ws.Range("K136").Value = 6639.900000000001
$ws.Range("M136").Value = -4089.900000000001
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 546.875
$ws.Range("I80").Value = 105.72727
$ws.Range("K80").Value = 105.72727
$ws.Range("M80").Value = 892.27273
$ws.Range("H83").Value = 546.875
$ws.Range("I83").Value = 105.72727
$ws.Range("K83").Value = 528.63635
$ws.Range("M83").Value = 4463.36365
$ws.Range("H107").Value = 1265.6086
$ws.Range("I107").Value = 1221.15
$ws.Range("J107").Value = 1562
$ws.Range("K107").Value = 1221.15
$ws.Range("L107").Value = 1562
$ws.Range("M107").Value = 698.8499999999999
$ws.Range("N107").Value = -5402
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 771661.5600000001
$ws.Range("I19").Value = 3333866.8
$ws.Range("J19").Value = 3000
$ws.Range("K19").Value = 3333866.8
$ws.Range("L19").Value = 3000
$ws.Range("M19").Value = -3333696.8
$ws.Range("N19").Value = -3340
$ws.Range("H24").Value = 771661.5600000001
$ws.Range("I24").Value = 3333866.8
$ws.Range("J24").Value = 3000
$ws.Range("K24").Value = 3333866.8
$ws.Range("L24").Value = 3000
$ws.Range("M24").Value = -3333696.8
$ws.Range("N24").Value = -3340
$ws.Range("H31").Value = 5237.885
$ws.Range("I31").Value = 3354.1316
$ws.Range("K31").Value = 3354.1316
$ws.Range("M31").Value = -3059.1316
$ws.Range("H34").Value = 5237.885
$ws.Range("I34").Value = 3354.1316
$ws.Range("K34").Value = 3354.1316
$ws.Range("M34").Value = -3152.1316
$ws.Range("H58").Value = 3235.5
$ws.Range("I58").Value = 3240.5
$ws.Range("K58").Value = 3240.5
$ws.Range("M58").Value = -3037.5
$ws.Range("H132").Value = 1437.6666
$ws.Range("I132").Value = 575.4
$ws.Range("K132").Value = 1726.2
$ws.Range("M132").Value = 803.8000000000002
$ws.Range("H136").Value = 3235.5
$ws.Range("I136").Value = 3240.5
$ws.Range("K136").Value = 9721.5
$ws.Range("M136").Value = -7171.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 377.42105
$ws.Range("J34").Value = 199.25
$ws.Range("L34").Value = 597.75
$ws.Range("N34").Value = -765.75
$ws.Range("H39").Value = 100
$ws.Range("I39").Value = 100
$ws.Range("K39").Value = 300
$ws.Range("M39").Value = -6
$ws.Range("H64").Value = 1000
$ws.Range("I64").Value = 1000
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 3000
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = -2730
$ws.Range("N64").Value = ""
$ws.Range("H67").Value = 1000
$ws.Range("I67").Value = 1000
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 3000
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = -2064
$ws.Range("N67").Value = ""
$ws.Range("H86").Value = 333638.28
$ws.Range("I86").Value = 1000197.4
$ws.Range("J86").Value = 358.7
$ws.Range("K86").Value = 3000592.2
$ws.Range("L86").Value = 1076.1
$ws.Range("M86").Value = -2999406.2
$ws.Range("N86").Value = -3448.1
$ws.Range("H89").Value = 333638.28
$ws.Range("I89").Value = 1000197.4
$ws.Range("J89").Value = 358.7
$ws.Range("K89").Value = 9001776.6
$ws.Range("L89").Value = 3228.3
$ws.Range("M89").Value = -8995848.6
$ws.Range("N89").Value = -15084.3
$ws.Range("H109").Value = 3841.25
$ws.Range("I109").Value = 2957.1428
$ws.Range("K109").Value = 8871.428400000001
$ws.Range("M109").Value = -7831.428400000001
$ws.Range("H116").Value = 34832.668
$ws.Range("I116").Value = 50999.5
$ws.Range("K116").Value = 152998.5
$ws.Range("M116").Value = -149556.5
$ws.Range("H132").Value = 9348.895
$ws.Range("I132").Value = 13648.5
$ws.Range("J132").Value = 1978.1428
$ws.Range("K132").Value = 122836.5
$ws.Range("L132").Value = 17803.2852
$ws.Range("M132").Value = -120306.5
$ws.Range("N132").Value = -22863.2852
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 25208.5
$ws.Range("I10").Value = 5313
$ws.Range("J10").Value = 64999.5
$ws.Range("K10").Value = 5313
$ws.Range("L10").Value = 64999.5
$ws.Range("M10").Value = -5144
$ws.Range("N10").Value = -65337.5
$ws.Range("H80").Value = 10209.75
$ws.Range("I80").Value = 7560.75
$ws.Range("K80").Value = 7560.75
$ws.Range("M80").Value = -6562.75
$ws.Range("H83").Value = 10209.75
$ws.Range("I83").Value = 7560.75
$ws.Range("K83").Value = 37803.75
$ws.Range("M83").Value = -32811.75
$ws.Range("H107").Value = 789.8461
$ws.Range("I107").Value = 434.7143
$ws.Range("J107").Value = 1204.1666
$ws.Range("K107").Value = 434.7143
$ws.Range("L107").Value = 1204.1666
$ws.Range("M107").Value = 1485.2857
$ws.Range("N107").Value = -5044.1666
$ws.Range("H127").Value = 70163
$ws.Range("I127").Value = 100000
$ws.Range("J127").Value = 40326
$ws.Range("K127").Value = 100000
$ws.Range("L127").Value = 40326
$ws.Range("M127").Value = -95040
$ws.Range("N127").Value = -50246
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").Value = ""
$ws.Range("H136").Value = 13195.5
$ws.Range("J136").Value = 13195.5
$ws.Range("L136").Value = 39586.5
$ws.Range("N136").Value = -44686.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1501
$ws.Range("I55").Value = 1501
$ws.Range("K55").Value = 1501
$ws.Range("M55").Value = -1328
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").Value = ""
$ws.Range("H136").Value = 4443.1875
$ws.Range("I136").Value = 3372.375
$ws.Range("K136").Value = 10117.125
$ws.Range("M136").Value = -7567.125
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 109166.664
$ws.Range("J46").Value = 109166.664
$ws.Range("L46").Value = 109166.664
$ws.Range("N46").Value = -109628.664
$ws.Range("H122").Value = 5244.1714
$ws.Range("I122").Value = 4895.839
$ws.Range("K122").Value = 14687.517
$ws.Range("M122").Value = -12237.517
$ws.Range("H126").Value = 4290.143
$ws.Range("I126").Value = 2006.4
$ws.Range("K126").Value = 6019.200000000001
$ws.Range("M126").Value = -3549.200000000001
$ws.Range("H134").Value = 109166.664
$ws.Range("J134").Value = 109166.664
$ws.Range("L134").Value = 327499.992
$ws.Range("N134").Value = -332569.992
$ws.Range("H136").Value = 1598.5
$ws.Range("I136").Value = 1598.5
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 4795.5
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -2245.5
$ws.Range("N136").Value = ""
